$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Iranzamin bank (row 24) - table_reading finished, Excel_File_Name filled in
# (written before Sina so the new shared-string entries land in the same
# order as the authoritative edit: iranzamin then sina)
$ws.Range("D24").Value = "iranzamin_bank_branches_20241117"
$ws.Range("A24:E24").ClearFormats()

# Sina bank (row 21) - table_reading finished, Excel_File_Name filled in
$ws.Range("D21").Value = "sina_bank_branches_20241117"
$ws.Range("A21:E21").ClearFormats()

# Khavarmianeh bank (row 25) - drop the highlight now that it is no longer the
# most-recently-completed row
$ws.Range("A25:E25").ClearFormats()

# Move the active selection to D15, matching where work continues next
$ws.Range("D15").Select()
